$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SearchInputData")

# Update B2 value from "Asus Laptop" to "Apple" (new shared string)
$ws.Range("B2").Value = "Apple"

# Update selection on the active sheet view
$ws.Range("C11").Select()
